$wb = $excel.ActiveWorkbook

# "OFF" sheet - Road ("R") row, update Week 15/16 cumulative stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 590
$wsOff.Range("C3").Value = 424
$wsOff.Range("D3").Value = 162
$wsOff.Range("E3").Value = 76
$wsOff.Range("F3").Value = 7

# "DEF" sheet - Road ("R") row, update Week 15/16 cumulative stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 600
$wsDef.Range("C3").Value = 431
$wsDef.Range("D3").Value = 106
$wsDef.Range("E3").Value = 56
